# Update the languages workbook: add a new "toggle" key/translation row,
# and move the view (frozen-pane scroll + active selection) further down
# the sheet, mirroring the author's manual edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new localization row (A31:C31) -------------------------------
$lastRow = 31

$ws.Cells.Item($lastRow, 1).Value = "toggle"
$ws.Cells.Item($lastRow, 2).Value = "Toggle"
$ws.Cells.Item($lastRow, 3).Value = "切換"

# Column C carries a distinct cell style (Microsoft JhengHei font) in every
# data row; copy it from the row above instead of touching Font directly so
# we reuse the existing style record rather than minting a new one.
$ws.Range("C30").Copy()
$ws.Range("C31").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Update the window/view state ------------------------------------------
# The author scrolled the frozen pane down and moved the active selection
# in the bottom-right pane to B34 (just past the new data).
$excel.Goto($ws.Range("B34"), $true)

$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 2
